$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Gangbanger Tuna Burger
$ws.Range("A3").Value = 45284
$ws.Range("B3").Value = "Gangbanger Tuna Burger"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 2312240001
$ws.Range("E3").Value = "Josefe Johnatan M. Gillego"

# Row 4: Hardcore Overload
$ws.Range("A4").Value = 45284
$ws.Range("B4").Value = "Hardcore Overload "
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2312240001
$ws.Range("E4").Value = "Josefe Johnatan M. Gillego"

# Row 5: Chicano Chili
$ws.Range("A5").Value = 45284
$ws.Range("B5").Value = "Chicano Chili"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 2312240002
$ws.Range("E5").Value = "Josefe Johnatan M. Gillego"

# Row 6: Rastaparay Veg
$ws.Range("A6").Value = 45284
$ws.Range("B6").Value = "Rastaparay Veg"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 2312240003
$ws.Range("E6").Value = "Josefe Johnatan M. Gillego"

# Row 7: Chicano Chili
$ws.Range("A7").Value = 45284
$ws.Range("B7").Value = "Chicano Chili"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 2312240003
$ws.Range("E7").Value = "Josefe Johnatan M. Gillego"

# Row 8: Hardcore Overload
$ws.Range("A8").Value = 45284
$ws.Range("B8").Value = "Hardcore Overload "
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 2312240004
$ws.Range("E8").Value = "Jerome"

# Rows 9-14 are cleared out (only empty numeric placeholder remains in column A)
$ws.Range("A9:E14").ClearContents()

# Update the active selection to B4
$ws.Range("B4").Select() | Out-Null
